# Auto-generated edit script: updates cryptos price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.843.16"
$ws.Range("E2").Value = "  -1.58%  "
$ws.Range("D3").Value = "2.285.69"
$ws.Range("E3").Value = "  -2.42%  "
$ws.Range("E4").Value = "  +0.01%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "316.41"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.63%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "102.61"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -5.00%  "
$ws.Range("E7").Value = "  -0.99%  "
$ws.Range("E8").Value = "  +0.08%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.602"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -2.21%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "38.97"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -4.96%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0904"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -2.11%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "8.31"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -1.91%  "
$ws.Range("E13").Value = "  +0.05%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.963"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -2.64%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "15.24"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -3.82%  "
$ws.Range("D16").Value = "2.629.59"
$ws.Range("E16").Value = "  -2.33%  "
$ws.Range("D17").Value = "2.276.62"
$ws.Range("E17").Value = "  -2.58%  "
$ws.Range("D18").Value = "41.808.62"
$ws.Range("E18").Value = "  -1.50%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "7.55"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.65%  "
$ws.Range("E20").Value = "  -0.28%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "284.92"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +11.27%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "73.69"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -3.10%  "
$ws.Range("E23").Value = "  -0.37%  "
$ws.Range("E24").Value = "  -1.53%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "9.91"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +5.81%  "
$ws.Range("E26").Value = "  +0.64%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "10.76"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -5.05%  "
$ws.Range("E28").Value = "  +3.20%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "23.06"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +1.27%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "162.92"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -5.74%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "34.69"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -5.44%  "
$ws.Range("E32").Value = "  -1.12%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "2.91"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +2.33%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "5.84"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -3.31%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.131"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.29%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.115"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -7.28%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "4.55"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -1.12%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "2.89"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +9.67%  "
$ws.Range("E39").Value = "  -3.63%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "3.60"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -7.70%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "102.81"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +20.47%  "
$ws.Range("E42").Value = "  +0.05%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "69.64"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -1.40%  "
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.225"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -4.11%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.33%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "115.43"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +3.56%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "11.96"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.35%  "
$ws.Range("B48").Value = "ordi"
$ws.Range("C48").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "76.82"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +3.21%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "9.00"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -1.60%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "5.30"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -2.39%  "
$ws.Range("E51").Value = "  -1.09%  "
